$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column K: ratio (J-I)/I for rows 3..21, formatted as percentage ---
$ws.Range("K3").Formula = "=(J3-I3)/I3"
$ws.Range("K4:K21").Formula = "=(J4-I4)/I4"
$ws.Range("K3:K21").NumberFormat = "0.00%"
$ws.Range("K3:K21").Style = "Percent"

# --- K24: array formula averaging the absolute ratios ---
$ws.Range("K24").FormulaArray = "=SUM(ABS(K3:K21))/COUNT(K3:K21)"
$ws.Range("K24").NumberFormat = "0.00%"

# --- Column K width ---
$ws.Range("K1").ColumnWidth = 11.85546875

# --- Selection / view tidy-up ---
$ws.Range("I3:K21").Select()

# --- Chart: title + series names + y-axis minimum ---
$chart = $ws.ChartObjects(1).Chart
$chart.HasTitle = $true
$chart.ChartTitle.Text = "Vamia Gateway with Normal Antenna - Obstacles"
$chart.SeriesCollection(1).Name = "Real RSSI"
$chart.SeriesCollection(2).Name = "Computed RSSI"
$chart.Axes(2).MinimumScale = 30
